$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content changes -------------------------------------------------

# Row 3: drop the trailing period from the security-issue label.
$ws.Range("B3").Value = "Unmasked password login field"

# Row 4: previously-blank row now documents "Not clearing login fields".
$ws.Range("B4").Value = "Not clearing login fields"
$ws.Range("C4").Value = "Login.java"
$ws.Range("D4").Value = "Allows unauthorized user/s to relogin with the previous credentials of a recent user"
$ws.Range("E4").Value = 'Sets the fields to an empty string after clicking login button: frame.loginPnl.usernameFld.setText("");
frame.loginPnl.passwordFld.setText("");'

# Row 5: previously-blank row now documents "No login verification".
$ws.Range("B5").Value = "No login verification"
$ws.Range("C5").Value = "Login.java, Main.java"
$ws.Range("D5").Value = "Allows unauthorized user/s to login without having the need to input valid user credentials"
$ws.Range("E5").Value = 'Added event listener to login button that checks if the inputted username and password matches an existing account stored in the database. 
Method implemented:
public boolean validateLogin(String username, String password)'

# --- Formatting changes ------------------------------------------------

# The whole data body now wraps text (previously only some rows did).
$ws.Range("A3:E22").WrapText = $true

# New rows need taller rows to fit the wrapped text.
$ws.Rows(4).RowHeight = 46.8
$ws.Rows(5).RowHeight = 93.6

# Column D widened (and best-fit) to accommodate the new, longer text.
$ws.Columns("D").ColumnWidth = 65.296875
$ws.Columns("D").AutoFit()
$ws.Columns("D").ColumnWidth = 65.296875

# Update the remembered selection to match the saved workbook state.
$ws.Range("F5").Select()
